$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$found = "Era presente nell'Excel"

$ws.Range("C8").Value = $found
$ws.Range("C9").Value = $found
$ws.Range("C36").Value = $found
$ws.Range("C37").Value = $found
